$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Scalar fixes to existing cells ---
$ws.Cells.Item(54, 17).Value = 0   # Q54: 2 -> 0
$ws.Cells.Item(69, 17).Value = 0   # Q69: 2 -> 0
$ws.Cells.Item(1485, 15).Value = 1 # O1485: 0 -> 1
$ws.Cells.Item(1487, 18).Value = 0 # R1487: inlineStr(empty) -> 0
$ws.Cells.Item(1488, 18).Value = 0 # R1488: inlineStr(empty) -> 0

# --- Append 19 new weekly rows (1489-1507), columns B..Q bulk-written via 2D array ---
$data = New-Object 'object[,]' 19,16
$data[0,0] = 305.7000122070312
$data[0,1] = 308.5
$data[0,2] = 301.9500122070312
$data[0,3] = 306.6499938964844
$data[0,4] = 297.13232421875
$data[0,5] = 28994214
$data[0,6] = 2024
$data[0,7] = 7
$data[0,8] = 1
$data[0,9] = 0
$data[0,10] = 0
$data[0,11] = 0
$data[0,12] = 27
$data[0,13] = 0
$data[0,14] = 0
$data[0,15] = 0
$data[1,0] = 304.6000061035156
$data[1,1] = 310.9500122070312
$data[1,2] = 293.3999938964844
$data[1,3] = 304.5499877929688
$data[1,4] = 295.0975036621094
$data[1,5] = 71389999
$data[1,6] = 2024
$data[1,7] = 7
$data[1,8] = 8
$data[1,9] = 0
$data[1,10] = 0
$data[1,11] = 0
$data[1,12] = 28
$data[1,13] = 2
$data[1,14] = 0
$data[1,15] = 0
$data[2,0] = 306
$data[2,1] = 319.25
$data[2,2] = 301.7999877929688
$data[2,3] = 303.7999877929688
$data[2,4] = 294.3707885742188
$data[2,5] = 76231453
$data[2,6] = 2024
$data[2,7] = 7
$data[2,8] = 15
$data[2,9] = 0
$data[2,10] = 0
$data[2,11] = 0
$data[2,12] = 29
$data[2,13] = 0
$data[2,14] = 2
$data[2,15] = 2
$data[3,0] = 304.5499877929688
$data[3,1] = 329.3999938964844
$data[3,2] = 296.1000061035156
$data[3,3] = 328.7999877929688
$data[3,4] = 318.5948486328125
$data[3,5] = 120717745
$data[3,6] = 2024
$data[3,7] = 7
$data[3,8] = 22
$data[3,9] = 0
$data[3,10] = 0
$data[3,11] = 0
$data[3,12] = 30
$data[3,13] = 0
$data[3,14] = 0
$data[3,15] = 0
$data[4,0] = 333.3999938964844
$data[4,1] = 359.0499877929688
$data[4,2] = 331.0499877929688
$data[4,3] = 347.1000061035156
$data[4,4] = 336.3268737792969
$data[4,5] = 115325192
$data[4,6] = 2024
$data[4,7] = 7
$data[4,8] = 29
$data[4,9] = 0
$data[4,10] = 0
$data[4,11] = 0
$data[4,12] = 31
$data[4,13] = 0
$data[4,14] = 0
$data[4,15] = 0
$data[5,0] = 342
$data[5,1] = 349.6499938964844
$data[5,2] = 329.7999877929688
$data[5,3] = 333.3999938964844
$data[5,4] = 323.0520629882812
$data[5,5] = 78551198
$data[5,6] = 2024
$data[5,7] = 8
$data[5,8] = 5
$data[5,9] = 0
$data[5,10] = 0
$data[5,11] = 0
$data[5,12] = 32
$data[5,13] = 0
$data[5,14] = 0
$data[5,15] = 0
$data[6,0] = 331.7999877929688
$data[6,1] = 336.7999877929688
$data[6,2] = 320.6499938964844
$data[6,3] = 332.5
$data[6,4] = 332.5
$data[6,5] = 38683065
$data[6,6] = 2024
$data[6,7] = 8
$data[6,8] = 12
$data[6,9] = 0
$data[6,10] = 0
$data[6,11] = 0
$data[6,12] = 33
$data[6,13] = 0
$data[6,14] = 0
$data[6,15] = 0
$data[7,0] = 337.7999877929688
$data[7,1] = 358
$data[7,2] = 336.3500061035156
$data[7,3] = 352.2000122070312
$data[7,4] = 352.2000122070312
$data[7,5] = 81676646
$data[7,6] = 2024
$data[7,7] = 8
$data[7,8] = 19
$data[7,9] = 0
$data[7,10] = 0
$data[7,11] = 0
$data[7,12] = 34
$data[7,13] = 0
$data[7,14] = 0
$data[7,15] = 1
$data[8,0] = 355.4500122070312
$data[8,1] = 365.5
$data[8,2] = 343.5
$data[8,3] = 357.6499938964844
$data[8,4] = 357.6499938964844
$data[8,5] = 57018464
$data[8,6] = 2024
$data[8,7] = 8
$data[8,8] = 26
$data[8,9] = 0
$data[8,10] = 0
$data[8,11] = 0
$data[8,12] = 35
$data[8,13] = 0
$data[8,14] = 0
$data[8,15] = 0
$data[9,0] = 361
$data[9,1] = 367.2000122070312
$data[9,2] = 351.0499877929688
$data[9,3] = 352.1499938964844
$data[9,4] = 352.1499938964844
$data[9,5] = 78022614
$data[9,6] = 2024
$data[9,7] = 9
$data[9,8] = 2
$data[9,9] = 0
$data[9,10] = 0
$data[9,11] = 0
$data[9,12] = 36
$data[9,13] = 0
$data[9,14] = 0
$data[9,15] = 0
$data[10,0] = 352.1000061035156
$data[10,1] = 353.5499877929688
$data[10,2] = 338.5499877929688
$data[10,3] = 342.2999877929688
$data[10,4] = 342.2999877929688
$data[10,5] = 56509071
$data[10,6] = 2024
$data[10,7] = 9
$data[10,8] = 9
$data[10,9] = 0
$data[10,10] = 0
$data[10,11] = 0
$data[10,12] = 37
$data[10,13] = 0
$data[10,14] = 0
$data[10,15] = 0
$data[11,0] = 344.8999938964844
$data[11,1] = 344.8999938964844
$data[11,2] = 322.9500122070312
$data[11,3] = 331.2000122070312
$data[11,4] = 331.2000122070312
$data[11,5] = 43391815
$data[11,6] = 2024
$data[11,7] = 9
$data[11,8] = 16
$data[11,9] = 0
$data[11,10] = 0
$data[11,11] = 0
$data[11,12] = 38
$data[11,13] = 0
$data[11,14] = 0
$data[11,15] = 0
$data[12,0] = 333
$data[12,1] = 370.5
$data[12,2] = 331.2999877929688
$data[12,3] = 367.2999877929688
$data[12,4] = 367.2999877929688
$data[12,5] = 69540176
$data[12,6] = 2024
$data[12,7] = 9
$data[12,8] = 23
$data[12,9] = 0
$data[12,10] = 0
$data[12,11] = 0
$data[12,12] = 39
$data[12,13] = 0
$data[12,14] = 0
$data[12,15] = 0
$data[13,0] = 367.1499938964844
$data[13,1] = 376
$data[13,2] = 334.3500061035156
$data[13,3] = 340.25
$data[13,4] = 340.25
$data[13,5] = 72241572
$data[13,6] = 2024
$data[13,7] = 9
$data[13,8] = 30
$data[13,9] = 0
$data[13,10] = 0
$data[13,11] = 0
$data[13,12] = 40
$data[13,13] = 1
$data[13,14] = 0
$data[13,15] = 0
$data[14,0] = 342
$data[14,1] = 350
$data[14,2] = 328.25
$data[14,3] = 337.6499938964844
$data[14,4] = 337.6499938964844
$data[14,5] = 40689958
$data[14,6] = 2024
$data[14,7] = 10
$data[14,8] = 7
$data[14,9] = 0
$data[14,10] = 0
$data[14,11] = 0
$data[14,12] = 41
$data[14,13] = 0
$data[14,14] = 0
$data[14,15] = 0
$data[15,0] = 340
$data[15,1] = 355.8999938964844
$data[15,2] = 334.2999877929688
$data[15,3] = 342.5
$data[15,4] = 342.5
$data[15,5] = 60576583
$data[15,6] = 2024
$data[15,7] = 10
$data[15,8] = 14
$data[15,9] = 0
$data[15,10] = 0
$data[15,11] = 0
$data[15,12] = 42
$data[15,13] = 0
$data[15,14] = 0
$data[15,15] = 0
$data[16,0] = 342.5499877929688
$data[16,1] = 345.5
$data[16,2] = 303.1000061035156
$data[16,3] = 306.2999877929688
$data[16,4] = 306.2999877929688
$data[16,5] = 38911903
$data[16,6] = 2024
$data[16,7] = 10
$data[16,8] = 21
$data[16,9] = 0
$data[16,10] = 0
$data[16,11] = 0
$data[16,12] = 43
$data[16,13] = 0
$data[16,14] = 0
$data[16,15] = 0
$data[17,0] = 307.2999877929688
$data[17,1] = 316.2000122070312
$data[17,2] = 302.25
$data[17,3] = 313
$data[17,4] = 313
$data[17,5] = 40071931
$data[17,6] = 2024
$data[17,7] = 10
$data[17,8] = 28
$data[17,9] = 0
$data[17,10] = 0
$data[17,11] = 0
$data[17,12] = 44
$data[17,13] = 0
$data[17,14] = 0
$data[17,15] = 0
$data[18,0] = 313
$data[18,1] = 319
$data[18,2] = 298.1000061035156
$data[18,3] = 310.4500122070312
$data[18,4] = 310.4500122070312
$data[18,5] = 51082717
$data[18,6] = 2024
$data[18,7] = 11
$data[18,8] = 4
$data[18,9] = 0
$data[18,10] = 0
$data[18,11] = 0
$data[18,12] = 45
$data[18,13] = 0
$data[18,14] = 0
$data[18,15] = 0

$ws.Range("B1489:Q1507").Value = $data

# --- Column A (Datetime, formatted) for the new rows ---
$ws.Cells.Item(1489, 1).Value = 45474
$ws.Cells.Item(1489, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1490, 1).Value = 45481
$ws.Cells.Item(1490, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1491, 1).Value = 45488
$ws.Cells.Item(1491, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1492, 1).Value = 45495
$ws.Cells.Item(1492, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1493, 1).Value = 45502
$ws.Cells.Item(1493, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1494, 1).Value = 45509
$ws.Cells.Item(1494, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1495, 1).Value = 45516
$ws.Cells.Item(1495, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1496, 1).Value = 45523
$ws.Cells.Item(1496, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1497, 1).Value = 45530
$ws.Cells.Item(1497, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1498, 1).Value = 45537
$ws.Cells.Item(1498, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1499, 1).Value = 45544
$ws.Cells.Item(1499, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1500, 1).Value = 45551
$ws.Cells.Item(1500, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1501, 1).Value = 45558
$ws.Cells.Item(1501, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1502, 1).Value = 45565
$ws.Cells.Item(1502, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1503, 1).Value = 45572
$ws.Cells.Item(1503, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1504, 1).Value = 45579
$ws.Cells.Item(1504, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1505, 1).Value = 45586
$ws.Cells.Item(1505, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1506, 1).Value = 45593
$ws.Cells.Item(1506, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1507, 1).Value = 45600
$ws.Cells.Item(1507, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
